$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.61
$summary.Range("B4").Value = -0.39
$summary.Range("B5").Value = -0.78
$summary.Range("B6").Value = 10
$summary.Range("B8").Value = 7
$summary.Range("B9").Value = 30

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.61
$status.Range("D4").Value = 10
$status.Range("E4").Value = -0.39
$status.Range("F4").Value = -0.39
$status.Range("G4").Value = 30

# --- New trade row (#10) to append on both "All Trades" and "MarketMaking" sheets ---
$newRow = @(10, "2026-02-17", "13:08:30", "MarketMaking", "UP", 0.08, 0.058113, "CLOSED", -27.3594, -0.02, 99.61, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    # Dates/times are stored as plain text in this sheet (matching rows 2-10),
    # so force text format on those two columns before writing, to avoid Excel
    # auto-converting "2026-02-17" / "13:08:30" into date/time serials.
    $ws.Range("B11").NumberFormat = "@"
    $ws.Range("C11").NumberFormat = "@"
    for ($col = 1; $col -le $newRow.Length; $col++) {
        $ws.Cells.Item(11, $col).Value = $newRow[$col - 1]
    }
    # Restore the default "Normal" style so the new cells don't carry a
    # leftover Text number format now that the literal strings are stored.
    $ws.Range("B11:C11").Style = "Normal"
}
